$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsredovisning")

# Move the existing "Summa" totals row (row 22) down to row 25,
# carrying its formatting (bold) along with it.
$ws.Range("A22:B22").Cut($ws.Range("A25:B25"))

# The cut left behind a formatted-but-empty row 22; clear that leftover
# formatting so the new data rows start from the default style.
$ws.Range("A22:B22").ClearFormats()

# Apply the date formatting (same as the rows above) to the new rows 22-24,
# copying the format only (not the value) from A21 so no new style gets
# added to the style table.
$ws.Range("A21").Copy()
$ws.Range("A22:A24").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 22: new entry
$ws.Cells.Item(22, 1).Value = 45317
$ws.Cells.Item(22, 2).Value = 1
$ws.Cells.Item(22, 3).Value = "Spara ny post tester"

# Row 23: new entry
$ws.Cells.Item(23, 1).Value = 45317
$ws.Cells.Item(23, 2).Value = 1
$ws.Cells.Item(23, 3).Value = "Kontrollera indata tester"

# Row 24: new entry
$ws.Cells.Item(24, 1).Value = 45317
$ws.Cells.Item(24, 2).Value = 2
$ws.Cells.Item(24, 3).Value = "Hämta enskild uppgift + Test"

# Update the SUBTOTAL formula on the (now relocated) totals row to include the new rows
$ws.Range("B25").Formula = "=SUBTOTAL(109,B2:B23)"

# Grow the Excel table to include the three new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C25"))

$ws.Range("C24").Select()
